{"js": "// Add a new \"referencedDistributionID\" row to the end of the first table\n// (the \"Objet error\" field table), right after the \"sourceMessage\" row.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\ntable.addRows(\"End\", 1, [\n  [\n    \"referencedDistributionID\",\n    \"DistributionID referenc\u00e9\",\n    \"string\",\n    \"1..1\",\n    \"DistributionID du message source\",\n    \"\",\n  ],\n]);\n\nawait context.sync();\n", "ps1": "# Add a new \"referencedDistributionID\" row to the end of the first table\n# (the \"Objet error\" field table), right after the \"sourceMessage\" row.\n$doc = $word.ActiveDocument\n$table = $doc.Tables.Item(1)\n\n$newRow = $table.Rows.Add()\n$rowIndex = $table.Rows.Count\n\n$table.Cell($rowIndex, 1).Range.Text = \"referencedDistributionID\"\n$table.Cell($rowIndex, 2).Range.Text = \"DistributionID referenc\u00e9\"\n$table.Cell($rowIndex, 3).Range.Text = \"string\"\n$table.Cell($rowIndex, 4).Range.Text = \"1..1\"\n$table.Cell($rowIndex, 5).Range.Text = \"DistributionID du message source\"\n"}
